# Add MAVE-NN paper to the selected publications worksheet, and
# tweak the "Algebraic and Diagrammatic Methods" description + header styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Prepend "WOW! " to the description of the "Algebraic and Diagrammatic
#    Methods for the Rule-Based Modeling of Multiparticle Complexes" row (row 2, col G).
$ws.Range("G2").Value = 'WOW! This paper introduces an operator algebra framework that bridges two previously disconnected approaches to modeling multiparticle complexes in stochastic chemical systems: the statistical physics formalism (rooted in Doi''s 1976 Fock space approach) and rule-based computational methods developed for simulating biochemical complexes. The new formalism extends Fock space to support not just particle creation and annihilation but also the assembly and disassembly of multiparticle complexes, with rules specified by algebraic operators via Wick''s theorem and aided by diagrammatic tools. The result is a unified mathematical and computational framework applicable to both equilibrium and nonequilibrium systems, including a stochastic simulation algorithm for the latter.'

# 2. Remove the hyperlink that currently sits on E6 (the PNAS / "deep
#    sequencing" row) before we shift rows around, so we can recreate it
#    cleanly at its new location afterwards.
$ws.Hyperlinks.Delete()

# 3. Insert a new row above row 6 (PNAS "deep sequencing" paper), pushing
#    that row - and everything below it - down by one to make room for the
#    new MAVE-NN entry.
$ws.Rows("6:6").Insert()

# 4. Populate the new row 6 with the MAVE-NN paper details.
$ws.Range("A6").Value = 'MAVE-NN: learning genotype-phenotype maps from multiplex assays of variant effect'
$ws.Range("B6").Value = 'Ammar Tareen, Mahdi Kooshkbaghi, Anna Posfai, William T Ireland, David M McCandlish, Justin B Kinney'
$ws.Range("C6").Value = 'Genome biology 23 (1), 98, 2022'
$ws.Range("D6").Value = 2022
$ws.Range("E6").Value = 'https://link.springer.com/article/10.1186/s13059-022-02661-7'
$ws.Range("G6").Value = 'hey'

# 5. Re-create the hyperlink on the PNAS row, which is now row 7.
$ws.Hyperlinks.Add($ws.Range("E7"), 'https://www.pnas.org/doi/10.1073/pnas.1004290107')
$ws.Range("E7").Style = "Hyperlink"

# 6. The header row no longer carries the bold header style.
$ws.Range("A1:G1").Font.Bold = $false

# 7. Restore the cursor/selection like the saved file shows.
$ws.Range("A12").Select()
